$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.127.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.125.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.54%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.83%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.114.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.63%  "
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.473"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.86%  "
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.641.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.120.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "63.001.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "456.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.739"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("E27").Value = "  +5.22%  "
$ws.Range("E28").Value = "  +7.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.61%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.114"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0807"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.22%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.69%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "430.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.53%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0374"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.01%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.937.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.95%  "
$ws.Range("E44").Value = "  +11.93%  "
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("E46").Value = "  +10.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.05%  "
